# Shop.xlsx fix: add the missing "Count" column (I) used by the shop plugin.
#
# Columns A:H already hold Id/Type/ItemID/Gold/Steel/Stone/Diamond/Level plus
# the metadata rows (Type/Public/Private/Save/Cache/Ref/Upload/Desc) in rows
# 1-9, followed by per-shop-entry data in rows 10-71. This adds a parallel
# "Count" field (an int column) with the same metadata-row semantics as the
# existing columns, defaulted to 1 for every shop entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry over the header/metadata-row formatting from column H (rows 2-8) so
# the new column matches the look of the rest of the table. Row 1's header
# cell intentionally keeps the plain row-default look (no copy needed there).
$ws.Range("H2:H8").Copy()
$ws.Range("I2:I8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header
$ws.Range("I1").Value = "Count"

# Metadata rows mirroring the other columns' field definitions
$ws.Range("I2").Value = "int"     # Type
$ws.Range("I3").Value = $false    # Public
$ws.Range("I4").Value = $false    # Private
$ws.Range("I5").Value = $true     # Save
$ws.Range("I6").Value = $false    # Cache
$ws.Range("I7").Value = $false    # Ref
$ws.Range("I8").Value = $false    # Upload
# Row 9 ("Desc") is left blank, matching the other columns.

# Data rows: every shop entry defaults to a count of 1.
for ($r = 10; $r -le 71; $r++) {
  $ws.Cells.Item($r, 9).Value = 1
}
